$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: force a text-typed cell write so Excel does not
# auto-coerce numeric-looking strings (e.g. "1.010") into numbers,
# which would silently drop meaningful trailing zeros. Switch the
# cell to Text format, assign the literal string, then clear the
# format back off so no stray styling is left behind.

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Update price (column D) and volume (column E) values per row ---
Set-TextValue $ws.Range("D2") "28.454.58"
$ws.Range("E2").Value = "  -3.81%  "

Set-TextValue $ws.Range("D3") "1.956.91"
$ws.Range("E3").Value = "  -2.39%  "

Set-TextValue $ws.Range("D4") "1.010"
$ws.Range("E4").Value = "  -0.43%  "

Set-TextValue $ws.Range("D5") "321.75"
$ws.Range("E5").Value = "  -2.55%  "

Set-TextValue $ws.Range("D6") "1.010"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("E7").Value = "  -4.78%  "

Set-TextValue $ws.Range("D8") "0.4062"
$ws.Range("E8").Value = "  -3.80%  "

Set-TextValue $ws.Range("D9") "53.21"
$ws.Range("E9").Value = "  -2.66%  "

Set-TextValue $ws.Range("D10") "0.08453"
$ws.Range("E10").Value = "  -6.30%  "

Set-TextValue $ws.Range("D11") "1.059"
$ws.Range("E11").Value = "  -5.21%  "

Set-TextValue $ws.Range("D12") "22.12"
$ws.Range("E12").Value = "  -4.99%  "

Set-TextValue $ws.Range("D13") "1.970.85"
$ws.Range("E13").Value = "  -5.08%  "

Set-TextValue $ws.Range("D14") "7.633"
$ws.Range("E14").Value = "  -5.37%  "

Set-TextValue $ws.Range("D15") "6.187"
$ws.Range("E15").Value = "  -4.33%  "

Set-TextValue $ws.Range("D16") "1.013"
$ws.Range("E16").Value = "  -0.18%  "

Set-TextValue $ws.Range("D17") "0.00001075"
$ws.Range("E17").Value = "  -3.61%  "

Set-TextValue $ws.Range("D18") "89.20"
$ws.Range("E18").Value = "  -5.59%  "

Set-TextValue $ws.Range("D19") "0.06613"
$ws.Range("E19").Value = "  -1.04%  "

Set-TextValue $ws.Range("D20") "18.70"
$ws.Range("E20").Value = "  -4.93%  "

Set-TextValue $ws.Range("D21") "1.010"
$ws.Range("E21").Value = "  -0.36%  "

Set-TextValue $ws.Range("D22") "5.824"
$ws.Range("E22").Value = "  -2.56%  "

Set-TextValue $ws.Range("D23") "28.485.47"
$ws.Range("E23").Value = "  -3.86%  "

Set-TextValue $ws.Range("D24") "11.59"
$ws.Range("E24").Value = "  -3.49%  "

Set-TextValue $ws.Range("D25") "2.291"
$ws.Range("E25").Value = "  -0.67%  "

Set-TextValue $ws.Range("D26") "2.186.22"
$ws.Range("E26").Value = "  -5.67%  "

Set-TextValue $ws.Range("D27") "154.39"
$ws.Range("E27").Value = "  -2.77%  "

Set-TextValue $ws.Range("D28") "20.23"
$ws.Range("E28").Value = "  -2.48%  "

Set-TextValue $ws.Range("D29") "5.991"
$ws.Range("E29").Value = "  -5.47%  "

Set-TextValue $ws.Range("D30") "2.164"
$ws.Range("E30").Value = "  -5.83%  "

Set-TextValue $ws.Range("D31") "123.80"
$ws.Range("E31").Value = "  -3.42%  "

Set-TextValue $ws.Range("D32") "0.9866"
$ws.Range("E32").Value = "  -6.66%  "

Set-TextValue $ws.Range("D33") "0.09605"
$ws.Range("E33").Value = "  -3.45%  "

Set-TextValue $ws.Range("D34") "1.448"
$ws.Range("E34").Value = "  -7.44%  "

Set-TextValue $ws.Range("D35") "5.594"
$ws.Range("E35").Value = "  -4.17%  "

Set-TextValue $ws.Range("D36") "3.662"
$ws.Range("E36").Value = "  -3.67%  "

Set-TextValue $ws.Range("D37") "0.02343"
$ws.Range("E37").Value = "  -4.92%  "

# --- Rows 38 and 39: FraxShare and Hedera swap places (data refreshed) ---
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D38") "8.808"
$ws.Range("E38").Value = "  -5.18%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D39") "0.06226"
$ws.Range("E39").Value = "  -2.95%  "

Set-TextValue $ws.Range("D40") "1.257"
$ws.Range("E40").Value = "  -3.64%  "

Set-TextValue $ws.Range("D41") "0.6232"
$ws.Range("E41").Value = "  -5.01%  "

$ws.Range("E42").Value = "  -4.52%  "

$ws.Range("E43").Value = "  -0.27%  "

Set-TextValue $ws.Range("D44") "0.1925"
$ws.Range("E44").Value = "  -5.93%  "

Set-TextValue $ws.Range("D45") "1.338"
$ws.Range("E45").Value = "  +2.86%  "

Set-TextValue $ws.Range("D46") "0.5972"
$ws.Range("E46").Value = "  -5.88%  "

Set-TextValue $ws.Range("D47") "12.99"
$ws.Range("E47").Value = "  -3.37%  "

$ws.Range("E48").Value = "  -6.16%  "

Set-TextValue $ws.Range("D49") "3.401"
$ws.Range("E49").Value = "  -3.16%  "

Set-TextValue $ws.Range("D50") "0.00000000329"
$ws.Range("E50").Value = "  -3.09%  "

Set-TextValue $ws.Range("D51") "0.06833"
$ws.Range("E51").Value = "  -2.26%  "
